$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-07 Tuesday" "2025-10-08 Wednesday"

Replace-Text "480÷9=" "450÷3="
Replace-Text "435÷2=" "437÷9="
Replace-Text "664÷6=" "595÷4="
Replace-Text "314÷7=" "221÷3="
Replace-Text "295÷5=" "679÷4="
Replace-Text "857÷3=" "955÷5="
Replace-Text "939÷3=" "113÷8="
Replace-Text "480÷8=" "917÷8="
Replace-Text "715÷5=" "106÷7="
Replace-Text "703÷9=" "519÷3="
Replace-Text "545÷7=" "958÷2="
Replace-Text "763÷4=" "433÷3="
Replace-Text "491÷3=" "753÷5="
Replace-Text "848÷3=" "500÷9="
Replace-Text "562÷8=" "766÷4="
Replace-Text "588÷8=" "937÷4="
Replace-Text "954÷4=" "175÷3="
Replace-Text "471÷6=" "769÷6="
Replace-Text "983÷4=" "119÷5="
Replace-Text "440÷4=" "975÷9="
Replace-Text "518÷5=" "913÷7="
Replace-Text "421÷3=" "713÷9="
Replace-Text "864÷3=" "915÷6="
Replace-Text "879÷6=" "333÷7="
Replace-Text "860÷6=" "292÷2="
